$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column price cells that look like plain numbers need an explicit Text
# format first, otherwise Excel auto-converts the assigned string into a
# number (and can drop formatting such as a trailing zero, e.g. "15.50").

$ws.Range("D2").Value = "25.723.18"
$ws.Range("E2").Value = "  -1.03%  "

$ws.Range("D3").Value = "1.625.17"
$ws.Range("E3").Value = "  -0.88%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.53"
$ws.Range("E5").Value = "  -0.01%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5095"
$ws.Range("E6").Value = "  +0.18%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2553"
$ws.Range("E8").Value = "  -0.11%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06285"
$ws.Range("E9").Value = "  -0.82%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.35"
$ws.Range("E10").Value = "  -0.90%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07766"
$ws.Range("E11").Value = "  -0.06%  "

$ws.Range("D12").Value = "1.632.43"
$ws.Range("E12").Value = "  -0.41%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.216"
$ws.Range("E13").Value = "  -1.27%  "

$ws.Range("D14").Value = "1.847.96"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5501"
$ws.Range("E15").Value = "  +1.71%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.41"
$ws.Range("E16").Value = "  -1.04%  "

$ws.Range("D17").Value = "0.0₅7465"
$ws.Range("E17").Value = "  -2.78%  "

$ws.Range("D18").Value = "25.742.23"
$ws.Range("E18").Value = "  -1.02%  "

$ws.Range("E19").Value = "  +0.09%  "

$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "193.55"
$ws.Range("E20").Value = "  -2.40%  "

$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.392"
$ws.Range("E21").Value = "  -0.37%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.745"
$ws.Range("E22").Value = "  -1.33%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.981"
$ws.Range("E23").Value = "  -0.79%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.002"
$ws.Range("E24").Value = "  -0.11%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.883"
$ws.Range("E25").Value = "  +1.13%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "141.52"
$ws.Range("E26").Value = "  +0.40%  "

$ws.Range("E27").Value = "  +5.71%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.50"
$ws.Range("E28").Value = "  -0.75%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.698"
$ws.Range("E29").Value = "  -1.43%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.235"
$ws.Range("E30").Value = "  +0.19%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.04856"
$ws.Range("E31").Value = "  -0.51%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.218"
$ws.Range("E32").Value = "  -0.82%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.143"
$ws.Range("E33").Value = "  -0.53%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.533"
$ws.Range("E34").Value = "  +1.01%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.372"
$ws.Range("E35").Value = "  +0.29%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.8906"
$ws.Range("E36").Value = "  -1.16%  "

$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5486"
$ws.Range("E37").Value = "  +1.20%  "

$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.531"
$ws.Range("E38").Value = "  -2.13%  "

$ws.Range("D39").Value = "1.108.83"
$ws.Range("E39").Value = "  -2.95%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01541"
$ws.Range("E40").Value = "  -0.99%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.001"
$ws.Range("E41").Value = "  -0.02%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.537"
$ws.Range("E42").Value = "  +2.52%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7949"
$ws.Range("E43").Value = "  -1.57%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "97.03"
$ws.Range("E44").Value = "  -2.10%  "

$ws.Range("D45").Value = "1.774.03"
$ws.Range("E45").Value = "  -0.26%  "

$ws.Range("D46").Value = "0.0₈111"
$ws.Range("E46").Value = "  -12.28%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4423"
$ws.Range("E47").Value = "  -2.33%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9976"
$ws.Range("E48").Value = "  -0.55%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "54.46"
$ws.Range("E49").Value = "  -0.63%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05123"
$ws.Range("E50").Value = "  +0.42%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.483"
$ws.Range("E51").Value = "  +2.30%  "
